$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("current")

# New column F header ("Equal") and values (0.5) for each state row (2-52)
$ws.Range("F1").Value = "Equal"
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 6).Value = 0.5
}

# Widen columns D and E to fit the new "Economist_prob_red" / "538_prob_red" headers
$ws.Columns.Item(4).ColumnWidth = 17.833333333333336
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666

# Move the active selection to C7
$ws.Range("C7").Select()
